$d = $word.ActiveDocument

$d.Content.Find.Execute("korisnci", $true, $false, $false, $false, $false,
                         $true, 1, $false, "korisnici", 2)
